# Update sheet name to reflect the new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-02-15"

# Update the label text for the February row
$ws.Range("A3").Value = "February (through 02-15)"

# Update February row (row 3) values
$ws.Range("C3").Value = 21
$ws.Range("D3").Value = 35
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 69
$ws.Range("I3").Value = 69

# Update Total row (row 4) values
$ws.Range("C4").Value = 72
$ws.Range("D4").Value = 110
$ws.Range("E4").Value = 116
$ws.Range("F4").Value = 62
$ws.Range("G4").Value = 111
$ws.Range("H4").Value = 286
$ws.Range("I4").Value = 230
